# Update "想去人数" (number of people who want to go) values in column F
# on both the "展览" sheet and the "全部类型" sheet, reflecting a refreshed
# data pull (output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" — rows keyed by their row number in that sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    3  = 534
    4  = 1100
    6  = 38
    8  = 39
    10 = 15797
    11 = 247
    14 = 6202
    21 = 29
    23 = 20
    24 = 13
    27 = 869
    29 = 5005
    31 = 11095
    32 = 1233
    35 = 179
    36 = 3811
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型" — same events, but shifted by one or two rows because this
# sheet also aggregates rows from the "演出" sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    4  = 534
    5  = 1100
    7  = 38
    9  = 39
    11 = 15797
    12 = 247
    15 = 6202
    22 = 29
    24 = 20
    25 = 13
    28 = 869
    30 = 5005
    33 = 11095
    34 = 1233
    37 = 179
    38 = 3811
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
